$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1014.38023562374
$ws.Range("C1").Value = 9.8443438331496207
$ws.Range("G1").Value = 9.8418471971860697
$ws.Range("Q1").Value = 9.8584060415702996
$ws.Range("T1").Value = 0.19700172088752901
$ws.Range("U1").Value = 9.9381617420716495
$ws.Range("AE1").Value = 9.8584060415702996
$ws.Range("AH1").Value = 8.18457624922371
$ws.Range("AI1").Value = 9.9381617420716495
$ws.Range("AL1").Value = 5.3002636771839402
$ws.Range("AP1").Value = 5.3002636771839402
$ws.Range("B2").Value = 2028.76
$ws.Range("C2").Value = 19.762205894692976
$ws.Range("G2").Value = 19.743004418991941
$ws.Range("Q2").Value = 23.8440548036751
$ws.Range("R2").Value = 3.5971410228619098
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 0.322335745314967
$ws.Range("U2").Value = 23.839561083151398
$ws.Range("V2").Value = 1
$ws.Range("AE2").Value = 20.3216216338122
$ws.Range("AH2").Value = 22.748767013319899
$ws.Range("AI2").Value = 20.329103308037499
$ws.Range("AL2").Value = 9.9033685267321747
$ws.Range("AP2").Value = 9.9033685267321747
$ws.Range("B3").Value = 3043.14
$ws.Range("C3").Value = 29.115972251621969
$ws.Range("G3").Value = 29.085727412886488
$ws.Range("Q3").Value = 44.481573839427902
$ws.Range("R3").Value = 14.400043945629999
$ws.Range("S3").Value = 4
$ws.Range("T3").Value = 0.29507097673404298
$ws.Range("U3").Value = 44.303823883042099
$ws.Range("V3").Value = 4
$ws.Range("AE3").Value = 36.298838111504097
$ws.Range("AH3").Value = 35.456943046366
$ws.Range("AI3").Value = 36.170866732494197
$ws.Range("AL3").Value = 12.407563508323898
$ws.Range("AP3").Value = 12.407563508323898
$ws.Range("B4").Value = 4057.52
$ws.Range("C4").Value = 41.924944702572205
$ws.Range("G4").Value = 41.832294235166387
$ws.Range("Q4").Value = 62.516100007510701
$ws.Range("R4").Value = 21.604941364373701
$ws.Range("S4").Value = 6
$ws.Range("T4").Value = 0.54955808429978503
$ws.Range("U4").Value = 62.475549684763401
$ws.Range("V4").Value = 6
$ws.Range("AE4").Value = 48.553249208808502
$ws.Range("AH4").Value = 35.476287623673997
$ws.Range("AI4").Value = 48.4655525513923
$ws.Range("AL4").Value = 16.860545565823742
$ws.Range("AP4").Value = 16.860545565823742
$ws.Range("B5").Value = 5071.8999999999996
$ws.Range("C5").Value = 52.814924323620332
$ws.Range("G5").Value = 52.654640661424608
$ws.Range("Q5").Value = 84.235267111039803
$ws.Range("R5").Value = 32.428611939589302
$ws.Range("S5").Value = 9
$ws.Range("T5").Value = 0.50446656437209703
$ws.Range("U5").Value = 84.076017589274201
$ws.Range("V5").Value = 9
$ws.Range("AE5").Value = 63.214879555727698
$ws.Range("AH5").Value = 79.513050458824296
$ws.Range("AI5").Value = 63.0463327910334
$ws.Range("AL5").Value = 20.319508511507696
$ws.Range("AP5").Value = 20.319508511507696
$ws.Range("B6").Value = 6086.28
$ws.Range("C6").Value = 66.441826425685477
$ws.Range("G6").Value = 65.586733890776046
$ws.Range("Q6").Value = 101.577891914391
$ws.Range("R6").Value = 39.649695238677303
$ws.Range("S6").Value = 11
$ws.Range("T6").Value = 0.79661290144165797
$ws.Range("U6").Value = 101.51578143586499
$ws.Range("V6").Value = 11
$ws.Range("AE6").Value = 74.829097008049203
$ws.Range("AH6").Value = 108.890661222194
$ws.Range("AI6").Value = 74.813136737609895
$ws.Range("AL6").Value = 23.661184528677854
$ws.Range("AP6").Value = 23.661184528677854
$ws.Range("B7").Value = 7100.66
$ws.Range("C7").Value = 76.751976545525295
$ws.Range("G7").Value = 75.011646641741478
$ws.Range("Q7").Value = 125.92730933983199
$ws.Range("R7").Value = 54.117372179038398
$ws.Range("S7").Value = 15
$ws.Range("T7").Value = 0.60009394560478402
$ws.Range("U7").Value = 125.510468124703
$ws.Range("V7").Value = 15
$ws.Range("AE7").Value = 83.219744779148201
$ws.Range("AH7").Value = 85.0579397073069
$ws.Range("AI7").Value = 83.085036094464201
$ws.Range("AL7").Value = 27.182327681147452
$ws.Range("AP7").Value = 27.182327681147452
$ws.Range("B8").Value = 8115.04
$ws.Range("C8").Value = 90.488391899777923
$ws.Range("G8").Value = 89.394662569256539
$ws.Range("Q8").Value = 144.59864419796901
$ws.Range("R8").Value = 61.3652445723501
$ws.Range("S8").Value = 17
$ws.Range("T8").Value = 0.68014484531765296
$ws.Range("U8").Value = 144.37639901947099
$ws.Range("V8").Value = 17
$ws.Range("AE8").Value = 90.607213425234789
$ws.Range("AH8").Value = 91.801755810958397
$ws.Range("AI8").Value = 90.722340009400796
$ws.Range("AL8").Value = 34.878759329679497
$ws.Range("AP8").Value = 34.878759329679497
$ws.Range("B9").Value = 9129.42
$ws.Range("C9").Value = 97.767615081483342
$ws.Range("G9").Value = 95.817868587472077
$ws.Range("Q9").Value = 181.920901003202
$ws.Range("R9").Value = 75.869776419360207
$ws.Range("S9").Value = 21
$ws.Range("T9").Value = 0.84824110756258198
$ws.Range("U9").Value = 181.248012778175
$ws.Range("V9").Value = 21
$ws.Range("AE9").Value = 129.67412142128936
$ws.Range("AH9").Value = 190.73703013370499
$ws.Range("AI9").Value = 129.41786725040535
$ws.Range("AL9").Value = 42.033180460678423
$ws.Range("AP9").Value = 42.033180460678423
$ws.Range("B10").Value = 10143.799999999999
$ws.Range("C10").Value = 109.5009830679254
$ws.Range("G10").Value = 108.29160753283814
$ws.Range("Q10").Value = 460.67677276646799
$ws.Range("R10").Value = 90.378294656707695
$ws.Range("S10").Value = 25
$ws.Range("T10").Value = 0.70165356906384702
$ws.Range("U10").Value = 460.48804689691099
$ws.Range("V10").Value = 25
$ws.Range("AE10").Value = 180.765961489662
$ws.Range("AH10").Value = 210.01795238087701
$ws.Range("AI10").Value = 180.89311609569799
$ws.Range("AL10").Value = 48.945727191344986
$ws.Range("AP10").Value = 48.945727191344986

$ws.Range("AI15").Select()
